# Update Name of Algo
# Apply updated imputed values for columns C and D as produced by the KNN algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.833
$ws.Range("D3").Value = -7.757
$ws.Range("D5").Value = -7.953999999999999
$ws.Range("C9").Value = -11.899
$ws.Range("D11").Value = -8.15
$ws.Range("D12").Value = -7.644999999999999
$ws.Range("C13").Value = -12.201
$ws.Range("C16").Value = -12.439
$ws.Range("C18").Value = -12.362
$ws.Range("C20").Value = -12.32
$ws.Range("D21").Value = -7.891
